$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-21 Thursday", "2025-08-22 Friday"),
    @("33×47=1551", "94×12=1128"),
    @("51×91=4641", "88×17=1496"),
    @("48×67=3216", "20×38=760"),
    @("93×97=9021", "46×60=2760"),
    @("52×12=624", "19×95=1805"),
    @("62×67=4154", "53×61=3233"),
    @("66×72=4752", "41×94=3854"),
    @("66×96=6336", "23×83=1909"),
    @("29×54=1566", "49×24=1176"),
    @("64×46=2944", "54×18=972"),
    @("93×24=2232", "40×36=1440"),
    @("19×74=1406", "17×97=1649"),
    @("61×11=671", "74×20=1480"),
    @("49×57=2793", "74×47=3478"),
    @("35×44=1540", "56×61=3416"),
    @("58×63=3654", "82×91=7462"),
    @("16×24=384", "77×76=5852"),
    @("57×61=3477", "92×95=8740"),
    @("50×94=4700", "25×52=1300"),
    @("37×44=1628", "68×90=6120"),
    @("40×26=1040", "67×68=4556"),
    @("89×95=8455", "84×26=2184"),
    @("31×27=837", "84×39=3276"),
    @("28×97=2716", "60×92=5520"),
    @("74×76=5624", "58×31=1798")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done replacing $($replacements.Count) items"
